$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting from existing header cell (G1) onto the new H1 header
# so it picks up the same bold/centered/bordered style used by other headers.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New "Save" data column values for rows 2-5
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 1
